# The deck's slide-master theme ("Integral" / "Red Violet" colour scheme,
# ppt/theme/theme1.xml) is swapped for the default "Office Theme" colour
# scheme (the palette that, before this edit, only the Notes Master used,
# ppt/theme/theme2.xml).
#
# PowerPoint's automation surface has no "replace this theme's XML" verb,
# so this is done the same way a user would from the ribbon: by recolouring
# the twelve theme colour slots (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink)
# on the slide master's theme colour scheme, in order, to the target
# palette's RGB values.

function HexToRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme palette, in clrScheme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeTheme = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000",
    "4472C4", "70AD47", "0563C1", "954F72"
)

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = HexToRgb $officeTheme[$i - 1]
}
